# "organizing documents for demo"
#
# 1) Slide 1: the MCU label text box (shape id=1138, "TextBox 1137") is
#    updated from "MSP430FR2311" to "MSP430FR2355" (the "MCU" line below it
#    is left untouched).
# 2) The cached "datetimeFigureOut" field shown in the Date placeholder of
#    the slide master and every slide layout is refreshed from 4/1/2025 to
#    4/7/2025 (PowerPoint re-stamps this cached text whenever the deck is
#    saved with the placeholder set to auto-update).

$p = $ppt.ActivePresentation

# --- 1) Update the MCU part label on slide 1 ------------------------------

$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Id -eq 1138) {
        $tr = $shp.TextFrame.TextRange
        $firstPara = $tr.Paragraphs(1, 1)
        $label = $firstPara.Text.TrimEnd([char]13)
        if ($label -eq "MSP430FR2311") {
            # Go through a disjoint placeholder first so the engine doesn't
            # keep a stale run around for the unchanged "MSP430FR2" prefix
            # (it otherwise splits the text into two <a:r> runs).
            $firstPara.Text = "__TMP__"
            $firstPara2 = $tr.Paragraphs(1, 1)
            $firstPara2.Text = "MSP430FR2355"
        }
        break
    }
}

# --- 2) Refresh the cached date field on the master + every layout --------

function Update-DatePlaceholder($shapes, [string]$newDate) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        $isDatePlaceholder = $false
        try {
            if ($shp.HasTextFrame -and $shp.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder -and $shp.TextFrame.HasText) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

Update-DatePlaceholder $p.SlideMaster.Shapes "4/7/2025"

$layouts = $p.SlideMaster.CustomLayouts
for ($k = 1; $k -le $layouts.Count; $k++) {
    $layout = $layouts.Item($k)
    Update-DatePlaceholder $layout.Shapes "4/7/2025"
}
